$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 4
$ws.Range("B4").Value = 'Count of Licenses'
$ws.Range("C4").Value = 'NA'
$ws.Range("F4").Value = 'Counts already at the council district level'

# Row 5
$ws.Range("B5").Value = 'Count of Licenses'
$ws.Range("C5").Value = 'NA'
$ws.Range("F5").Value = 'Counts already at the council district level'

# Row 6
$ws.Range("B6").Value = 'Count of crashes'
$ws.Range("C6").Value = 'NA'
$ws.Range("F6").Value = 'Counts were summarized per council district using a spatial join'
$ws.Rows.Item(6).RowHeight = 62.4

# Row 11
$ws.Range("B11").Value = 'Percentage of buildings with an open housing code violation'
$ws.Range("C11").Value = 'The number of parcels of land available per council district'
$ws.Range("F11").Value = 'Code violations were set at the council district level. However land parcels had to be spatially joined to council districts. '
$ws.Rows.Item(11).RowHeight = 93.6

# Row 12
$ws.Range("B12").Value = 'Percentage of homes that lack a complete kitchen'
$ws.Range("C12").Value = 'Total Number of Occupied Housing Units'
$ws.Range("F12").Value = 'Census Tract aggregated up to the council district using Block population as weights'
$ws.Range("G12").Value = 'Variables Used:' + [char]10 + '#B25052_001 --> Total Occupied housing units tract' + [char]10 + '#B25052_003 --> Lacking complete kitchen facilities tract'
$ws.Rows.Item(12).RowHeight = 78

# Row 13
$ws.Range("B13").Value = 'Percentage of homes that lack complete plumbing'
$ws.Range("C13").Value = 'Total Number of Occupied Housing Units'
$ws.Range("F13").Value = 'Census Tract aggregated up to the council district using Block population as weights'
$ws.Range("G13").Value = 'Variables Used:' + [char]10 + '#B25048_001 --> Total Occupied housing units tract' + [char]10 + '#B25048_003 --> Lacking plumbing occupied housing units tract'
$ws.Rows.Item(13).RowHeight = 78

# Row 14
$ws.Range("B14").Value = 'Median Household Income'
$ws.Range("C14").Value = 'NA'
$ws.Range("F14").Value = 'Census Block Group aggregated up to the council district using total number of households as weights'
$ws.Range("G14").Value = 'Variables Used:' + [char]10 + '#B19013_001 - Estimate!!Median household income in the past 12 months (in 2022 inflation-adjusted dollars)' + [char]10 + ''
$ws.Rows.Item(14).RowHeight = 93.6

# Row 15
$ws.Range("B15").Value = 'Percentage of Households that Own their home'
$ws.Range("C15").Value = 'Total Number of Households'
$ws.Range("F15").Value = 'Census Block Group aggregated up to the council district using total number of households as weights'
$ws.Range("G15").Value = 'Variables Used:' + [char]10 + '#B25003_001 --> Estimate!!Total: block group' + [char]10 + '#B25003_002 --> Estimate!!Total:!!Owner occupied block group'
$ws.Rows.Item(15).RowHeight = 93.6

# Row 16
$ws.Range("B16").Value = 'Percentage of Households that Rent their home'
$ws.Range("C16").Value = 'Total Number of Households'
$ws.Range("F16").Value = 'Census Block Group aggregated up to the council district using total number of households as weights'
$ws.Range("G16").Value = 'Variables Used:' + [char]10 + '#B25003_001 --> Estimate!!Total: block group' + [char]10 + '#B25003_003 --> Estimate!!Total:!!Renter occupied block group '
$ws.Rows.Item(16).RowHeight = 93.6

# Row 25
$ws.Range("B25").Value = 'Total Number of Fatal Shootings per Council District'
$ws.Range("C25").Value = 'NA'
$ws.Range("F25").Value = 'Counts were summarized per council district using a spatial join'
$ws.Rows.Item(25).RowHeight = 62.4

# Row 26
$ws.Range("B26").Value = 'Total Number of Nonfatal Shootings per Council District'
$ws.Range("C26").Value = 'NA'
$ws.Range("F26").Value = 'Counts were summarized per council district using a spatial join'
$ws.Rows.Item(26).RowHeight = 62.4

# Row 27
$ws.Range("B27").Value = 'Percentage of people Uninsured'
$ws.Range("C27").Value = 'Total Number of People'
$ws.Range("F27").Value = 'Block Group variables were aggregated up to the council district level using Block population as the weights'
$ws.Range("G27").Value = 'Variables Used:' + [char]10 + '#B18135_007 -- Estimate!!Total:!!Under 19 years:!!With a disability:!!No health insurance coverage' + [char]10 + '#B18135_012 -- Estimate!!Total:!!Under 19 years:!!No disability:!!No health insurance coverage' + [char]10 + '#B18135_018 -- Estimate!!Total:!!19 to 64 years:!!With a disability:!!No health insurance coverage' + [char]10 + '#B18135_023 -- Estimate!!Total:!!19 to 64 years:!!No disability:!!No health insurance coverage' + [char]10 + '#B18135_029 -- Estimate!!Total:!!65 years and over:!!With a disability:!!No health insurance coverage' + [char]10 + '#B18135_034 -- Estimate!!Total:!!65 years and over:!!No disability:!!No health insurance coverage'
$ws.Rows.Item(27).RowHeight = 409.6

# Row 28
$ws.Range("B28").Value = 'Percentage of Council district covered by greenspace'
$ws.Range("C28").Value = 'Shape Area of Council District'
$ws.Range("F28").Value = '2010 Census tracts, converted to 2020 Blocks, then aggregated up to the council district level using 2010 Shape Area as weights'
$ws.Range("G28").Value = 'Used IPUMS crosswalk to convert 2010 Census tracts to 2020 Blocks. Before joining at the council district level'
$ws.Rows.Item(28).RowHeight = 109.2

# Row 29
$ws.Range("B29").Value = 'Heat Vulnerability Index'
$ws.Range("C29").Value = 'NA'
$ws.Range("F29").Value = '2010 Census tracts, converted to 2020 Blocks, then aggregated up to the council district level using 2010 Shape Area as weights'
$ws.Range("G29").Value = 'Used IPUMS crosswalk to convert 2010 Census tracts to 2020 Blocks. Before joining at the council district level'
$ws.Rows.Item(29).RowHeight = 109.2

# Update selection to match target state
[void]$ws.Range("G30").Select()
